$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.252
$ws.Range("G2").Value = 0.1003921568627451
$ws.Range("H2").Value = 0.1003921568627451
$ws.Range("I2").Value = -0.3176470588235294
$ws.Range("J2").Value = -0.3176470588235294
$ws.Range("K2").Value = -4.54
$ws.Range("L2").Value = -0.5934640522875817
$ws.Range("U2").Value = 1.06
$ws.Range("V2").Value = 0.1284848484848485
$ws.Range("W2").Value = -14.05572755417957
$ws.Range("X2").Value = 0.1817389212981077
$ws.Range("Y2").Value = -14.23746647547767
$ws.Range("Z2").Value = 0.2354063452010955
$ws.Range("AA2").Value = -0.07477613318152446
$ws.Range("AB2").Value = 0.07348724733250468
$ws.Range("AC2").Value = -0.1482633805140292
$ws.Range("AD2").Value = 23.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 23.4
$ws.Range("AG2").Value = 22.34
$ws.Range("AH2").Value = 0.7393364928909952
$ws.Range("AI2").Value = 0.9983361064891846
$ws.Range("AJ2").Value = 0.7303040209218699
$ws.Range("AK2").Value = 0.9982572947852897
$ws.Range("AL2").Value = 0.728
$ws.Range("AM2").Value = 0.726
$ws.Range("AN2").Value = -25.51799345692475
$ws.Range("AO2").Value = -3.337912087912088
$ws.Range("AP2").Value = -24.36205016357688
$ws.Range("AQ2").Value = -3.347107438016529

# Row 3
$ws.Range("D3").Value = -0.252
$ws.Range("G3").Value = 0.1003921568627451
$ws.Range("H3").Value = 0.1003921568627451
$ws.Range("I3").Value = -0.3176470588235294
$ws.Range("J3").Value = -0.3176470588235294
$ws.Range("K3").Value = -4.54
$ws.Range("L3").Value = -0.5934640522875817
$ws.Range("U3").Value = 1.06
$ws.Range("V3").Value = 0.1284848484848485
$ws.Range("W3").Value = -14.05572755417957
$ws.Range("X3").Value = 0.1817389212981077
$ws.Range("Y3").Value = -14.23746647547767
$ws.Range("Z3").Value = 0.2354063452010955
$ws.Range("AA3").Value = -0.07477613318152446
$ws.Range("AB3").Value = 0.07348724733250468
$ws.Range("AC3").Value = -0.1482633805140292
$ws.Range("AD3").Value = 23.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 23.4
$ws.Range("AG3").Value = 22.34
$ws.Range("AH3").Value = 0.7393364928909952
$ws.Range("AI3").Value = 0.9983361064891846
$ws.Range("AJ3").Value = 0.7303040209218699
$ws.Range("AK3").Value = 0.9982572947852897
$ws.Range("AL3").Value = 0.728
$ws.Range("AM3").Value = 0.726
$ws.Range("AN3").Value = -25.51799345692475
$ws.Range("AO3").Value = -3.337912087912088
$ws.Range("AP3").Value = -24.36205016357688
$ws.Range("AQ3").Value = -3.347107438016529

